$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.697.08"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "'2.243.84"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'112.68"
$ws.Range("E5").Value = "  -1.59%  "
$ws.Range("D6").Value = "'296.21"
$ws.Range("E6").Value = "  +7.15%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("D10").Value = "'43.96"
$ws.Range("E10").Value = "  -5.36%  "
$ws.Range("D11").Value = "'0.0923"
$ws.Range("E11").Value = "  -0.49%  "
$ws.Range("D12").Value = "'54.30"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").Value = "'9.01"
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("D14").Value = "'1.07"
$ws.Range("E14").Value = "  +22.11%  "
$ws.Range("E15").Value = "  -0.90%  "
$ws.Range("D16").Value = "'15.18"
$ws.Range("E16").Value = "  -0.71%  "
$ws.Range("D17").Value = "'2.583.94"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").Value = "'2.273.21"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("D19").Value = "'42.720.31"
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("D21").Value = "'7.18"
$ws.Range("E21").Value = "  +5.78%  "
$ws.Range("D22").Value = "'74.59"
$ws.Range("E22").Value = "  +3.28%  "
$ws.Range("D23").Value = "'3.48"
$ws.Range("E23").Value = "  +16.58%  "
$ws.Range("E24").Value = "  +2.98%  "
$ws.Range("D25").Value = "'250.34"
$ws.Range("E25").Value = "  +8.01%  "
$ws.Range("E26").Value = "  -3.22%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.63%  "
$ws.Range("D28").Value = "'11.54"
$ws.Range("E28").Value = "  -4.18%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.23"
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "'176.20"
$ws.Range("E30").Value = "  +1.51%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "'37.53"
$ws.Range("E31").Value = "  -7.03%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'22.06"
$ws.Range("E32").Value = "  +4.59%  "
$ws.Range("D33").Value = "'3.17"
$ws.Range("E33").Value = "  -3.31%  "
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("D35").Value = "'5.68"
$ws.Range("E35").Value = "  +1.77%  "
$ws.Range("D36").Value = "'5.09"
$ws.Range("E36").Value = "  +9.62%  "
$ws.Range("D37").Value = "'4.25"
$ws.Range("E37").Value = "  -3.99%  "
$ws.Range("D38").Value = "'0.128"
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("D39").Value = "'0.0377"
$ws.Range("E39").Value = "  +1.28%  "
$ws.Range("E40").Value = "  -1.89%  "
$ws.Range("E41").Value = "  -5.94%  "
$ws.Range("D42").Value = "'71.97"
$ws.Range("E42").Value = "  +1.45%  "
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").Value = "'12.45"
$ws.Range("E45").Value = "  -5.82%  "
$ws.Range("E46").Value = "  -1.02%  "
$ws.Range("D47").Value = "'5.50"
$ws.Range("E47").Value = "  -2.84%  "
$ws.Range("E48").Value = "  +2.66%  "
$ws.Range("D49").Value = "'105.61"
$ws.Range("E49").Value = "  +4.69%  "
$ws.Range("D50").Value = "'8.56"
$ws.Range("E50").Value = "  +1.40%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.0984"
$ws.Range("E51").Value = "  -0.45%  "
